# Update the ligand-receptor table with recomputed TPM-based statistics.
#
# The sheet originally held a single "Sending cluster" (MuSCs) x 4
# "Target cluster" rows (ECs, FAPs, MuSCs, Resolving-Mac). The refreshed
# scripts add a second sending cluster (ECs) against the same four target
# clusters, so the table grows from 4 data rows to 8: the original four
# rows are recalculated for the new TPM values, and four new rows are
# appended for the ECs sending cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E..T the numeric NATMI statistics.
$rows = @(
    @("ECs",   "a", "Mc1r", "ECs",           1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 3, 1, 0.7564453333333333, 2.269336, 0.2228868660598179, 0.2228868660598179, 0.034296474968,      0.308668274712,      0.03878885949634809, 0.03878885949634809),
    @("ECs",   "a", "Mc1r", "FAPs",          1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 3, 1, 1.780736333333333,  5.342209, 0.5246945458259835, 0.5246945458259834, 0.080736804617,      0.7266312415530001,  0.09131225799137997, 0.09131225799137996),
    @("ECs",   "a", "Mc1r", "MuSCs",         1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 3, 1, 0.6826836666666667, 2.048051, 0.2011529667359422, 0.2011529667359422, 0.030952194763,      0.278569752867,      0.03500652282445402, 0.03500652282445402),
    @("ECs",   "a", "Mc1r", "Resolving-Mac", 1, 0.3333333333333333, 0.045339, 0.136017, 0.1740293637846656, 0.1740293637846656, 3, 1, 0.173988,            0.5219640000000001, 0.05126562137825639, 0.05126562137825638, 0.007888441932000001, 0.07099597738800001, 0.00892172347248351, 0.00892172347248351),
    @("MuSCs", "a", "Mc1r", "ECs",           1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 3, 1, 0.7564453333333333, 2.269336, 0.2228868660598179, 0.2228868660598179, 0.1627764454986667,  1.464988009488,      0.1840980065634699,  0.1840980065634699),
    @("MuSCs", "a", "Mc1r", "FAPs",          1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 3, 1, 1.780736333333333,  5.342209, 0.5246945458259835, 0.5246945458259834, 0.3831895286246667,  3.448705757622,      0.4333822878346035,  0.4333822878346035),
    @("MuSCs", "a", "Mc1r", "MuSCs",         1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 3, 1, 0.6826836666666667, 2.048051, 0.2011529667359422, 0.2011529667359422, 0.1469039674953333,  1.322135707458,      0.1661464439114882,  0.1661464439114882),
    @("MuSCs", "a", "Mc1r", "Resolving-Mac", 1, 0.3333333333333333, 0.215186, 0.645558, 0.8259706362153344, 0.8259706362153345, 3, 1, 0.173988,            0.5219640000000001, 0.05126562137825639, 0.05126562137825638, 0.037439781768,      0.336958035912,      0.04234389790577288, 0.04234389790577288)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
